# Fruta / hortaliza, semanal
# Insert two new weekly records (date 2022-06-10 / serial 44722) above the
# existing block that starts at row 32, pushing the old rows 32-37 down to
# rows 34-39.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows at position 32 (existing rows 32:37 shift down to 34:39)
$ws.Rows("32:33").Insert()

# Common / shared values for the two new rows
$mercadoId = 1
$mercado   = "Agrícola del Norte S.A. de Arica"
$region    = "Arica y Parinacota"
$fecha     = 44722
$codreg    = 15
$tipo      = "Fruta"
$productoId = 100108
$producto  = "Tropicales y subtropicales"
$categoriaId = 100108001
$categoria = "Guayaba"
$variedad  = "Sin especificar"
$unidad    = "$/kilo (en caja de 10 kilos )"
$origen    = "Región de Arica y Parinacota"
$kgUnidad  = 1

# Row 32 - Calidad "Primera"
$ws.Cells.Item(32, 1).Value = $mercadoId
$ws.Cells.Item(32, 2).Value = $mercado
$ws.Cells.Item(32, 3).Value = $region
$ws.Cells.Item(32, 4).Value = $fecha
$ws.Cells.Item(32, 5).Value = $codreg
$ws.Cells.Item(32, 6).Value = $tipo
$ws.Cells.Item(32, 7).Value = $productoId
$ws.Cells.Item(32, 8).Value = $producto
$ws.Cells.Item(32, 9).Value = $categoriaId
$ws.Cells.Item(32, 10).Value = $categoria
$ws.Cells.Item(32, 11).Value = $variedad
$ws.Cells.Item(32, 12).Value = "Primera"
$ws.Cells.Item(32, 13).Value = 140
$ws.Cells.Item(32, 14).Value = 800
$ws.Cells.Item(32, 15).Value = 900
$ws.Cells.Item(32, 16).Value = 850
$ws.Cells.Item(32, 17).Value = $unidad
$ws.Cells.Item(32, 18).Value = $origen
$ws.Cells.Item(32, 19).Value = 850
$ws.Cells.Item(32, 20).Value = $kgUnidad

# Row 33 - Calidad "Segunda"
$ws.Cells.Item(33, 1).Value = $mercadoId
$ws.Cells.Item(33, 2).Value = $mercado
$ws.Cells.Item(33, 3).Value = $region
$ws.Cells.Item(33, 4).Value = $fecha
$ws.Cells.Item(33, 5).Value = $codreg
$ws.Cells.Item(33, 6).Value = $tipo
$ws.Cells.Item(33, 7).Value = $productoId
$ws.Cells.Item(33, 8).Value = $producto
$ws.Cells.Item(33, 9).Value = $categoriaId
$ws.Cells.Item(33, 10).Value = $categoria
$ws.Cells.Item(33, 11).Value = $variedad
$ws.Cells.Item(33, 12).Value = "Segunda"
$ws.Cells.Item(33, 13).Value = 200
$ws.Cells.Item(33, 14).Value = 700
$ws.Cells.Item(33, 15).Value = 800
$ws.Cells.Item(33, 16).Value = 750
$ws.Cells.Item(33, 17).Value = $unidad
$ws.Cells.Item(33, 18).Value = $origen
$ws.Cells.Item(33, 19).Value = 750
$ws.Cells.Item(33, 20).Value = $kgUnidad
